$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 18149.725  # H132: 23003.025 -> 18149.725
$ws.Cells.Item(132, 9).Value = 2670.9023  # I132: 3386.8386 -> 2670.9023
$ws.Cells.Item(132, 10).Value = 81612.89999999999  # J132: 90569.89 -> 81612.89999999999
$ws.Cells.Item(132, 11).Value = 8012.706900000001  # K132: 10160.5158 -> 8012.706900000001
$ws.Cells.Item(132, 12).Value = 244838.7  # L132: 271709.67 -> 244838.7
$ws.Cells.Item(132, 13).Value = -5482.706900000001  # M132: -7630.515800000001 -> -5482.706900000001
$ws.Cells.Item(132, 14).Value = -249898.7  # N132: -276769.67 -> -249898.7

$ws.Cells.Item(138, 8).Value = 2919.125  # H138: 3327.7576 -> 2919.125
$ws.Cells.Item(138, 9).Value = 1580.68  # I138: 1886.7894 -> 1580.68
$ws.Cells.Item(138, 10).Value = 5149.8667  # J138: 5283.357 -> 5149.8667
$ws.Cells.Item(138, 11).Value = 4742.04  # K138: 5660.3682 -> 4742.04
$ws.Cells.Item(138, 12).Value = 15449.6001  # L138: 15850.071 -> 15449.6001
$ws.Cells.Item(138, 13).Value = 397.96  # M138: -520.3681999999999 -> 397.96
$ws.Cells.Item(138, 14).Value = -25729.6001  # N138: -26130.071 -> -25729.6001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(31, 8).Value = 9948.75  # H31: 10428.363 -> 9948.75
$ws.Cells.Item(31, 9).Value = 2131.3333  # I31: 2452 -> 2131.3333
$ws.Cells.Item(31, 10).Value = 19999.715  # J31: 20000 -> 19999.715
$ws.Cells.Item(31, 11).Value = 2131.3333  # K31: 2452 -> 2131.3333
$ws.Cells.Item(31, 12).Value = 19999.715  # L31: 20000 -> 19999.715
$ws.Cells.Item(31, 13).Value = -1837.3333  # M31: -2158 -> -1837.3333
$ws.Cells.Item(31, 14).Value = -20587.715  # N31: -20588 -> -20587.715

$ws.Cells.Item(32, 8).Value = 3697.7314  # H32: 5375.511 -> 3697.7314
$ws.Cells.Item(32, 9).Value = 1838  # I32: 2681.282 -> 1838
$ws.Cells.Item(32, 10).Value = 19638.285  # J32: 22888 -> 19638.285
$ws.Cells.Item(32, 11).Value = 1838  # K32: 2681.282 -> 1838
$ws.Cells.Item(32, 12).Value = 19638.285  # L32: 22888 -> 19638.285
$ws.Cells.Item(32, 13).Value = -1551  # M32: -2394.282 -> -1551
$ws.Cells.Item(32, 14).Value = -20212.285  # N32: -23462 -> -20212.285

$ws.Cells.Item(45, 8).Value = 59084.43  # H45: 72749.94 -> 59084.43
$ws.Cells.Item(45, 9).Value = 72379.94  # I45: 81895.664 -> 72379.94
$ws.Cells.Item(45, 10).Value = 2578.5  # J45: 4157 -> 2578.5
$ws.Cells.Item(45, 11).Value = 72379.94  # K45: 81895.664 -> 72379.94
$ws.Cells.Item(45, 12).Value = 2578.5  # L45: 4157 -> 2578.5
$ws.Cells.Item(45, 13).Value = -72002.94  # M45: -81518.664 -> -72002.94
$ws.Cells.Item(45, 14).Value = -3332.5  # N45: -4911 -> -3332.5

$ws.Cells.Item(92, 8).Value = 0  # H92: 35000 -> 0
$ws.Cells.Item(92, 10).Value = 0  # J92: 35000 -> 0
$ws.Cells.Item(92, 12).Value = 0  # L92: 35000 -> 0
$ws.Cells.Item(92, 14).ClearContents()  # remove N92

$ws.Cells.Item(122, 8).Value = 1882.8334  # H122: 1701.8928 -> 1882.8334
$ws.Cells.Item(122, 9).Value = 1837.8462  # I122: 1582.7 -> 1837.8462
$ws.Cells.Item(122, 10).Value = 1999.8  # J122: 1999.875 -> 1999.8
$ws.Cells.Item(122, 11).Value = 5513.5386  # K122: 4748.1 -> 5513.5386
$ws.Cells.Item(122, 12).Value = 5999.4  # L122: 5999.625 -> 5999.4
$ws.Cells.Item(122, 13).Value = -3063.5386  # M122: -2298.1 -> -3063.5386
$ws.Cells.Item(122, 14).Value = -10899.4  # N122: -10899.625 -> -10899.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 3535.75  # H105: 3391.1667 -> 3535.75
$ws.Cells.Item(105, 9).Value = 2616.25  # I105: 2128.625 -> 2616.25
$ws.Cells.Item(105, 10).Value = 4455.25  # J105: 4401.2 -> 4455.25
$ws.Cells.Item(105, 11).Value = 2616.25  # K105: 2128.625 -> 2616.25
$ws.Cells.Item(105, 12).Value = 4455.25  # L105: 4401.2 -> 4455.25
$ws.Cells.Item(105, 13).Value = -869.25  # M105: -381.625 -> -869.25
$ws.Cells.Item(105, 14).Value = -7949.25  # N105: -7895.2 -> -7949.25

$ws.Cells.Item(134, 8).Value = 2693.7942  # H134: 2913.7932 -> 2693.7942
$ws.Cells.Item(134, 9).Value = 2509.4194  # I134: 2768.12 -> 2509.4194
$ws.Cells.Item(134, 10).Value = 4599  # J134: 3824.25 -> 4599
$ws.Cells.Item(134, 11).Value = 7528.2582  # K134: 8304.360000000001 -> 7528.2582
$ws.Cells.Item(134, 12).Value = 13797  # L134: 11472.75 -> 13797
$ws.Cells.Item(134, 13).Value = -4993.2582  # M134: -5769.360000000001 -> -4993.2582
$ws.Cells.Item(134, 14).Value = -18867  # N134: -16542.75 -> -18867

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(15, 8).Value = 15000  # H15: 10000 -> 15000
$ws.Cells.Item(15, 9).Value = 0  # I15: 10000 -> 0
$ws.Cells.Item(15, 10).Value = 15000  # J15: 0 -> 15000
$ws.Cells.Item(15, 11).Value = 0  # K15: 10000 -> 0
$ws.Cells.Item(15, 12).Value = 15000  # L15: 0 -> 15000
$ws.Cells.Item(15, 13).ClearContents()  # remove M15
$ws.Cells.Item(15, 14).Value = -15340  # add N15

$ws.Cells.Item(16, 8).Value = 1495.0769  # H16: 1790.5 -> 1495.0769
$ws.Cells.Item(16, 9).Value = 1582  # I16: 1970.3334 -> 1582
$ws.Cells.Item(16, 10).Value = 1440.75  # J16: 1682.6 -> 1440.75
$ws.Cells.Item(16, 11).Value = 1582  # K16: 1970.3334 -> 1582
$ws.Cells.Item(16, 12).Value = 1440.75  # L16: 1682.6 -> 1440.75
$ws.Cells.Item(16, 13).Value = -1295  # M16: -1683.3334 -> -1295
$ws.Cells.Item(16, 14).Value = -2014.75  # N16: -2256.6 -> -2014.75

$ws.Cells.Item(31, 8).Value = 5753848.5  # H31: 6673428 -> 5753848.5
$ws.Cells.Item(31, 9).Value = 3010.111  # I31: 2732.7273 -> 3010.111
$ws.Cells.Item(31, 10).Value = 6810125  # J31: 8554906 -> 6810125
$ws.Cells.Item(31, 11).Value = 3010.111  # K31: 2732.7273 -> 3010.111
$ws.Cells.Item(31, 12).Value = 6810125  # L31: 8554906 -> 6810125
$ws.Cells.Item(31, 13).Value = -2715.111  # M31: -2437.7273 -> -2715.111
$ws.Cells.Item(31, 14).Value = -6810715  # N31: -8555496 -> -6810715

$ws.Cells.Item(34, 8).Value = 5753848.5  # H34: 6673428 -> 5753848.5
$ws.Cells.Item(34, 9).Value = 3010.111  # I34: 2732.7273 -> 3010.111
$ws.Cells.Item(34, 10).Value = 6810125  # J34: 8554906 -> 6810125
$ws.Cells.Item(34, 11).Value = 3010.111  # K34: 2732.7273 -> 3010.111
$ws.Cells.Item(34, 12).Value = 6810125  # L34: 8554906 -> 6810125
$ws.Cells.Item(34, 13).Value = -2808.111  # M34: -2530.7273 -> -2808.111
$ws.Cells.Item(34, 14).Value = -6810529  # N34: -8555310 -> -6810529

$ws.Cells.Item(58, 8).Value = 1839.3846  # H58: 1877.1316 -> 1839.3846
$ws.Cells.Item(58, 9).Value = 1205.2273  # I58: 1279.7368 -> 1205.2273
$ws.Cells.Item(58, 10).Value = 2660.0588  # J58: 2474.5264 -> 2660.0588
$ws.Cells.Item(58, 11).Value = 1205.2273  # K58: 1279.7368 -> 1205.2273
$ws.Cells.Item(58, 12).Value = 2660.0588  # L58: 2474.5264 -> 2660.0588
$ws.Cells.Item(58, 13).Value = -1002.2273  # M58: -1076.7368 -> -1002.2273
$ws.Cells.Item(58, 14).Value = -3066.0588  # N58: -2880.5264 -> -3066.0588

$ws.Cells.Item(113, 8).Value = 1495.0769  # H113: 1790.5 -> 1495.0769
$ws.Cells.Item(113, 9).Value = 1582  # I113: 1970.3334 -> 1582
$ws.Cells.Item(113, 10).Value = 1440.75  # J113: 1682.6 -> 1440.75
$ws.Cells.Item(113, 11).Value = 1582  # K113: 1970.3334 -> 1582
$ws.Cells.Item(113, 12).Value = 1440.75  # L113: 1682.6 -> 1440.75
$ws.Cells.Item(113, 13).Value = 588  # M113: 199.6666 -> 588
$ws.Cells.Item(113, 14).Value = -5780.75  # N113: -6022.6 -> -5780.75

$ws.Cells.Item(122, 8).Value = 241365.2  # H122: 89097.08 -> 241365.2
$ws.Cells.Item(122, 9).Value = 301228  # I122: 147652.2 -> 301228
$ws.Cells.Item(122, 10).Value = 1914  # J122: 1264.4 -> 1914
$ws.Cells.Item(122, 11).Value = 903684  # K122: 442956.6 -> 903684
$ws.Cells.Item(122, 12).Value = 5742  # L122: 3793.2 -> 5742
$ws.Cells.Item(122, 13).Value = -901234  # M122: -440506.6 -> -901234
$ws.Cells.Item(122, 14).Value = -10642  # N122: -8693.200000000001 -> -10642

$ws.Cells.Item(134, 8).Value = 54269.742  # H134: 56322.42 -> 54269.742
$ws.Cells.Item(134, 9).Value = 1615.2941  # I134: 1660 -> 1615.2941
$ws.Cells.Item(134, 11).Value = 4845.8823  # K134: 4980 -> 4845.8823
$ws.Cells.Item(134, 13).Value = -2310.8823  # M134: -2445 -> -2310.8823

$ws.Cells.Item(136, 8).Value = 1839.3846  # H136: 1877.1316 -> 1839.3846
$ws.Cells.Item(136, 9).Value = 1205.2273  # I136: 1279.7368 -> 1205.2273
$ws.Cells.Item(136, 10).Value = 2660.0588  # J136: 2474.5264 -> 2660.0588
$ws.Cells.Item(136, 11).Value = 3615.6819  # K136: 3839.2104 -> 3615.6819
$ws.Cells.Item(136, 12).Value = 7980.176399999999  # L136: 7423.5792 -> 7980.176399999999
$ws.Cells.Item(136, 13).Value = -1065.6819  # M136: -1289.2104 -> -1065.6819
$ws.Cells.Item(136, 14).Value = -13080.1764  # N136: -12523.5792 -> -13080.1764

$ws.Cells.Item(138, 8).Value = 43599.75  # H138: 27233.334 -> 43599.75
$ws.Cells.Item(138, 10).Value = 43599.75  # J138: 27233.334 -> 43599.75
$ws.Cells.Item(138, 12).Value = 43599.75  # L138: 27233.334 -> 43599.75
$ws.Cells.Item(138, 14).Value = -53879.75  # N138: -37513.334 -> -53879.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 295.61765  # H12: 300.42426 -> 295.61765
$ws.Cells.Item(12, 10).Value = 448.33334  # J12: 466.64706 -> 448.33334
$ws.Cells.Item(12, 12).Value = 1345.00002  # L12: 1399.94118 -> 1345.00002
$ws.Cells.Item(12, 14).Value = -1691.00002  # N12: -1745.94118 -> -1691.00002

$ws.Cells.Item(34, 8).Value = 1234.3  # H34: 1147.375 -> 1234.3
$ws.Cells.Item(34, 9).Value = 0  # I34: 700 -> 0
$ws.Cells.Item(34, 10).Value = 1234.3  # J34: 1166.826 -> 1234.3
$ws.Cells.Item(34, 11).Value = 0  # K34: 2100 -> 0
$ws.Cells.Item(34, 12).Value = 3702.9  # L34: 3500.478 -> 3702.9
$ws.Cells.Item(34, 13).ClearContents()  # remove M34
$ws.Cells.Item(34, 14).Value = -3870.9  # N34: -3668.478 -> -3870.9

$ws.Cells.Item(102, 8).Value = 7799.6  # H102: 7699.8 -> 7799.6
$ws.Cells.Item(102, 10).Value = 7799.6  # J102: 7699.8 -> 7799.6
$ws.Cells.Item(102, 12).Value = 23398.8  # L102: 23099.4 -> 23398.8
$ws.Cells.Item(102, 14).Value = -28266.8  # N102: -27967.4 -> -28266.8

$ws.Cells.Item(133, 8).Value = 4905.5654  # H133: 4326 -> 4905.5654
$ws.Cells.Item(133, 9).Value = 2702.9  # I133: 3757.5 -> 2702.9
$ws.Cells.Item(133, 10).Value = 6599.923  # J133: 6600 -> 6599.923
$ws.Cells.Item(133, 11).Value = 8108.700000000001  # K133: 11272.5 -> 8108.700000000001
$ws.Cells.Item(133, 12).Value = 19799.769  # L133: 19800 -> 19799.769
$ws.Cells.Item(133, 13).Value = -3048.700000000001  # M133: -6212.5 -> -3048.700000000001
$ws.Cells.Item(133, 14).Value = -29919.769  # N133: -29920 -> -29919.769

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(99, 8).Value = 11538.875  # H99: 11173.174 -> 11538.875
$ws.Cells.Item(99, 10).Value = 19970.834  # J99: 19972.727 -> 19970.834
$ws.Cells.Item(99, 12).Value = 19970.834  # L99: 19972.727 -> 19970.834
$ws.Cells.Item(99, 14).Value = -24462.834  # N99: -24464.727 -> -24462.834

$ws.Cells.Item(102, 8).Value = 815.1429000000001  # H102: 843.0526 -> 815.1429000000001
$ws.Cells.Item(102, 10).Value = 807.4286  # J102: 910.4 -> 807.4286
$ws.Cells.Item(102, 12).Value = 807.4286  # L102: 910.4 -> 807.4286
$ws.Cells.Item(102, 14).Value = -4051.4286  # N102: -4154.4 -> -4051.4286

$ws.Cells.Item(122, 8).Value = 1599.5  # H122: 1622.2222 -> 1599.5
$ws.Cells.Item(122, 9).Value = 1599.375  # I122: 1628.5714 -> 1599.375
$ws.Cells.Item(122, 11).Value = 4798.125  # K122: 4885.7142 -> 4798.125
$ws.Cells.Item(122, 13).Value = -2348.125  # M122: -2435.7142 -> -2348.125

$ws.Cells.Item(132, 8).Value = 4532.9165  # H132: 3818.4375 -> 4532.9165
$ws.Cells.Item(132, 9).Value = 2842  # I132: 2323.3333 -> 2842
$ws.Cells.Item(132, 11).Value = 8526  # K132: 6969.999899999999 -> 8526
$ws.Cells.Item(132, 13).Value = -5996  # M132: -4439.999899999999 -> -5996

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 3388.9333  # H68: 3491.8147 -> 3388.9333
$ws.Cells.Item(68, 9).Value = 3326.1365  # I68: 3419.3 -> 3326.1365
$ws.Cells.Item(68, 10).Value = 3561.625  # J68: 3699 -> 3561.625
$ws.Cells.Item(68, 11).Value = 3326.1365  # K68: 3419.3 -> 3326.1365
$ws.Cells.Item(68, 12).Value = 3561.625  # L68: 3699 -> 3561.625
$ws.Cells.Item(68, 13).Value = -2577.1365  # M68: -2670.3 -> -2577.1365
$ws.Cells.Item(68, 14).Value = -5059.625  # N68: -5197 -> -5059.625

$ws.Cells.Item(71, 8).Value = 3388.9333  # H71: 3491.8147 -> 3388.9333
$ws.Cells.Item(71, 9).Value = 3326.1365  # I71: 3419.3 -> 3326.1365
$ws.Cells.Item(71, 10).Value = 3561.625  # J71: 3699 -> 3561.625
$ws.Cells.Item(71, 11).Value = 16630.6825  # K71: 17096.5 -> 16630.6825
$ws.Cells.Item(71, 12).Value = 17808.125  # L71: 18495 -> 17808.125
$ws.Cells.Item(71, 13).Value = -12886.6825  # M71: -13352.5 -> -12886.6825
$ws.Cells.Item(71, 14).Value = -25296.125  # N71: -25983 -> -25296.125

$ws.Cells.Item(82, 8).Value = 9260059  # H82: 9260037 -> 9260059
$ws.Cells.Item(82, 9).Value = 798  # I82: 0 -> 798
$ws.Cells.Item(82, 10).Value = 10417467  # J82: 9260037 -> 10417467
$ws.Cells.Item(82, 11).Value = 798  # K82: 0 -> 798
$ws.Cells.Item(82, 12).Value = 10417467  # L82: 9260037 -> 10417467
$ws.Cells.Item(82, 13).Value = -437  # add M82
$ws.Cells.Item(82, 14).Value = -10418189  # N82: -9260759 -> -10418189

$ws.Cells.Item(85, 8).Value = 9260059  # H85: 9260037 -> 9260059
$ws.Cells.Item(85, 9).Value = 798  # I85: 0 -> 798
$ws.Cells.Item(85, 10).Value = 10417467  # J85: 9260037 -> 10417467
$ws.Cells.Item(85, 11).Value = 798  # K85: 0 -> 798
$ws.Cells.Item(85, 12).Value = 10417467  # L85: 9260037 -> 10417467
$ws.Cells.Item(85, 13).Value = 450  # add M85
$ws.Cells.Item(85, 14).Value = -10419963  # N85: -9262533 -> -10419963

$ws.Cells.Item(93, 8).Value = 941.8148  # H93: 1012.88 -> 941.8148
$ws.Cells.Item(93, 9).Value = 785.1  # I93: 860.8889 -> 785.1
$ws.Cells.Item(93, 10).Value = 1034  # J93: 1098.375 -> 1034
$ws.Cells.Item(93, 11).Value = 785.1  # K93: 860.8889 -> 785.1
$ws.Cells.Item(93, 12).Value = 1034  # L93: 1098.375 -> 1034
$ws.Cells.Item(93, 13).Value = 462.9  # M93: 387.1111 -> 462.9
$ws.Cells.Item(93, 14).Value = -3530  # N93: -3594.375 -> -3530

$ws.Cells.Item(122, 8).Value = 79261.69500000001  # H122: 113533.664 -> 79261.69500000001
$ws.Cells.Item(122, 9).Value = 113344.11  # I122: 145114.72 -> 113344.11
$ws.Cells.Item(122, 10).Value = 2576.25  # J122: 3000 -> 2576.25
$ws.Cells.Item(122, 11).Value = 340032.33  # K122: 435344.16 -> 340032.33
$ws.Cells.Item(122, 12).Value = 7728.75  # L122: 9000 -> 7728.75
$ws.Cells.Item(122, 13).Value = -337582.33  # M122: -432894.16 -> -337582.33
$ws.Cells.Item(122, 14).Value = -12628.75  # N122: -13900 -> -12628.75

$ws.Cells.Item(132, 8).Value = 4299.885  # H132: 4238.3335 -> 4299.885
$ws.Cells.Item(132, 9).Value = 2845.4546  # I132: 2874.1428 -> 2845.4546
$ws.Cells.Item(132, 10).Value = 5366.467  # J132: 5707.4614 -> 5366.467
$ws.Cells.Item(132, 11).Value = 8536.363799999999  # K132: 8622.428400000001 -> 8536.363799999999
$ws.Cells.Item(132, 12).Value = 16099.401  # L132: 17122.3842 -> 16099.401
$ws.Cells.Item(132, 13).Value = -6006.363799999999  # M132: -6092.428400000001 -> -6006.363799999999
$ws.Cells.Item(132, 14).Value = -21159.401  # N132: -22182.3842 -> -21159.401

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 2446.1538  # H62: 2440 -> 2446.1538
$ws.Cells.Item(62, 10).Value = 2446.1538  # J62: 2440 -> 2446.1538
$ws.Cells.Item(62, 12).Value = 2446.1538  # L62: 2440 -> 2446.1538
$ws.Cells.Item(62, 14).Value = -3694.1538  # N62: -3688 -> -3694.1538

$ws.Cells.Item(65, 8).Value = 2446.1538  # H65: 2440 -> 2446.1538
$ws.Cells.Item(65, 10).Value = 2446.1538  # J65: 2440 -> 2446.1538
$ws.Cells.Item(65, 12).Value = 12230.769  # L65: 12200 -> 12230.769
$ws.Cells.Item(65, 14).Value = -18470.769  # N65: -18440 -> -18470.769

$ws.Cells.Item(104, 8).Value = 27310  # H104: 43010 -> 27310
$ws.Cells.Item(104, 10).Value = 27310  # J104: 43010 -> 27310
$ws.Cells.Item(104, 12).Value = 27310  # L104: 43010 -> 27310
$ws.Cells.Item(104, 14).Value = -34298  # N104: -49998 -> -34298

$ws.Cells.Item(122, 8).Value = 1933.0588  # H122: 2113.9333 -> 1933.0588
$ws.Cells.Item(122, 9).Value = 2108.2307  # I122: 2323.0908 -> 2108.2307
$ws.Cells.Item(122, 10).Value = 1363.75  # J122: 1538.75 -> 1363.75
$ws.Cells.Item(122, 11).Value = 6324.6921  # K122: 6969.2724 -> 6324.6921
$ws.Cells.Item(122, 12).Value = 4091.25  # L122: 4616.25 -> 4091.25
$ws.Cells.Item(122, 13).Value = -3874.6921  # M122: -4519.2724 -> -3874.6921
$ws.Cells.Item(122, 14).Value = -8991.25  # N122: -9516.25 -> -8991.25
